# Update column F (dSF) values on the active worksheet to match re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 9
    3  = -6
    4  = -1
    5  = -3
    6  = 3
    7  = 3
    8  = 3
    9  = -2
    10 = -2
    11 = 3
    13 = 1
    17 = -4
    18 = -6
    20 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
